# Update the "sequence to use" column (B1:B20), the fitness/penalty
# improvement values (D1, D2), and the last generation fit value (B21)
# to reflect a new best individual found by the genetic algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sequence = @(11, 15, 14, 8, 6, 5, 10, 19, 18, 0, 12, 7, 4, 1, 3, 17, 9, 2, 16, 13)

for ($i = 0; $i -lt $sequence.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $sequence[$i]
}

$ws.Range("D1").Value = 110.9395668781889
$ws.Range("D2").Value = 73.37153010941981
$ws.Range("B21").Value = 0.7863468976227408
